$d = $word.ActiveDocument

# Locate the paragraph containing the footer line
# "Ver no Jupiter Salvar em pdf Salvar em docx" - this anchors the block of
# boilerplate site-footer paragraphs that must be removed.
$verParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Ver no Jupiter*") {
        $verParaIndex = $i
        break
    }
}

if ($verParaIndex -gt 0) {
    # The paragraph right before it is the blank spacer paragraph that was added
    # together with the footer block, and the paragraph right after it is the
    # copyright/Jekyll notice line - both need to disappear along with it.
    $startPara = $d.Paragraphs.Item($verParaIndex - 1)
    $endPara   = $d.Paragraphs.Item($verParaIndex + 1)

    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
